# Carga 201511 Noviembre de 2015
# Rename trailing headers and append new "Tipo de Servicio / Tipo de Cobro /
# Precio / kilos Integrados / Kilo Excedido" columns to the report header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing header cells shift two columns to the right (BA/BB keep their
#     previous contents, now the last two header cells get brand-new labels) ---
$ws.Range("BA1").Value = "Activo"
$ws.Range("BB1").Value = "Fecha Activo"
$ws.Range("BC1").Value = "Tipo de Servicio"
$ws.Range("BD1").Value = "Tipo de Cobro"
$ws.Range("BE1").Value = "Precio"
$ws.Range("BF1").Value = "kilos Integrados"
$ws.Range("BG1").Value = "Kilo Excedido"

# --- Match the bold-white-on-blue header style used by every other header
#     cell on row 1 (copy formatting only, values already set above) ---
$ws.Range("A1").Copy()
$ws.Range("BC1:BG1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths (best-fit) for the (now shifted/new) trailing columns ---
$ws.Columns.Item(53).ColumnWidth = 8.307291666666666
$ws.Columns.Item(54).ColumnWidth = 16.451822916666668
$ws.Columns.Item(55).ColumnWidth = 20.166666666666668
$ws.Columns.Item(56).ColumnWidth = 17.736979166666668
$ws.Columns.Item(57).ColumnWidth = 8.166666666666666
$ws.Columns.Item(58).ColumnWidth = 20.022135416666668
$ws.Columns.Item(59).ColumnWidth = 17.022135416666668
